$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell as TEXT (not a number),
# without leaving behind any NumberFormat/quote-prefix style change on the cell.
# Trick: put ="<text>" as a formula (forces string result), then Copy +
# PasteSpecial(xlPasteValues = -4163) over itself to collapse it down to a
# plain string value/constant while the cell keeps its original style.
function Set-TextValue($cell, $text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$ws.Range("D2").Value = '23.413.73'
$ws.Range("E2").Value = '  +1.55%  '

$ws.Range("D3").Value = '1.626.35'
$ws.Range("E3").Value = '  +2.21%  '

Set-TextValue $ws.Range("D4") '0.9951'
$ws.Range("E4").Value = '  -0.70%  '

Set-TextValue $ws.Range("D5") '307.61'
$ws.Range("E5").Value = '  +1.83%  '

Set-TextValue $ws.Range("D6") '0.9979'
$ws.Range("E6").Value = '  -0.42%  '

Set-TextValue $ws.Range("D7") '0.3784'
$ws.Range("E7").Value = '  +0.29%  '

Set-TextValue $ws.Range("D8") '52.93'
$ws.Range("E8").Value = '  +3.94%  '

Set-TextValue $ws.Range("D9") '0.3665'
$ws.Range("E9").Value = '  +1.68%  '

$ws.Range("E10").Value = '  +3.56%  '

Set-TextValue $ws.Range("D11") '0.08186'
$ws.Range("E11").Value = '  +1.40%  '

Set-TextValue $ws.Range("D12") '0.9995'
$ws.Range("E12").Value = '  -0.28%  '

Set-TextValue $ws.Range("D13") '23.17'
$ws.Range("E13").Value = '  +4.80%  '

Set-TextValue $ws.Range("D14") '6.658'
$ws.Range("E14").Value = '  +2.45%  '

Set-TextValue $ws.Range("D15") '7.436'
$ws.Range("E15").Value = '  +2.08%  '

Set-TextValue $ws.Range("D16") '0.00001256'
$ws.Range("E16").Value = '  +2.18%  '

$ws.Range("D17").Value = '1.619.68'
$ws.Range("E17").Value = '  +1.81%  '

Set-TextValue $ws.Range("D18") '94.84'
$ws.Range("E18").Value = '  +2.34%  '

Set-TextValue $ws.Range("D19") '0.06953'
$ws.Range("E19").Value = '  +1.90%  '

Set-TextValue $ws.Range("D20") '18.35'
$ws.Range("E20").Value = '  +1.89%  '

Set-TextValue $ws.Range("D21") '6.583'
$ws.Range("E21").Value = '  +1.70%  '

Set-TextValue $ws.Range("D22") '0.9980'
$ws.Range("E22").Value = '  -0.51%  '

Set-TextValue $ws.Range("D23") '12.98'
$ws.Range("E23").Value = '  +0.81%  '

$ws.Range("D24").Value = '23.417.05'
$ws.Range("E24").Value = '  +1.54%  '

Set-TextValue $ws.Range("D25") '3.149'
$ws.Range("E25").Value = '  +11.43%  '

Set-TextValue $ws.Range("D26") '2.423'
$ws.Range("E26").Value = '  +2.13%  '

Set-TextValue $ws.Range("D27") '21.37'
$ws.Range("E27").Value = '  +2.02%  '

Set-TextValue $ws.Range("D28") '151.10'
$ws.Range("E28").Value = '  +1.68%  '

Set-TextValue $ws.Range("D29") '5.285'
$ws.Range("E29").Value = '  +1.31%  '

Set-TextValue $ws.Range("D30") '136.49'
$ws.Range("E30").Value = '  +2.17%  '

Set-TextValue $ws.Range("D31") '2.413'
$ws.Range("E31").Value = '  +2.56%  '

Set-TextValue $ws.Range("D32") '6.880'
$ws.Range("E32").Value = '  +4.61%  '

$ws.Range("D33").Value = '1.798.45'
$ws.Range("E33").Value = '  +1.78%  '

Set-TextValue $ws.Range("D34") '0.9706'
$ws.Range("E34").Value = '  +2.38%  '

Set-TextValue $ws.Range("D35") '0.02796'
$ws.Range("E35").Value = '  +4.18%  '

Set-TextValue $ws.Range("D36") '10.49'
$ws.Range("E36").Value = '  +3.46%  '

Set-TextValue $ws.Range("D37") '0.07439'
$ws.Range("E37").Value = '  +0.20%  '

Set-TextValue $ws.Range("D38") '6.237'
$ws.Range("E38").Value = '  +2.52%  '

Set-TextValue $ws.Range("D39") '0.2533'
$ws.Range("E39").Value = '  +1.97%  '

Set-TextValue $ws.Range("D40") '0.08848'
$ws.Range("E40").Value = '  +0.69%  '

Set-TextValue $ws.Range("D41") '1.407'
$ws.Range("E41").Value = '  +4.29%  '

Set-TextValue $ws.Range("D42") '0.7171'
$ws.Range("E42").Value = '  +3.24%  '

Set-TextValue $ws.Range("D43") '12.72'
$ws.Range("E43").Value = '  +4.86%  '

Set-TextValue $ws.Range("D44") '16.14'
$ws.Range("E44").Value = '  +8.01%  '

Set-TextValue $ws.Range("D45") '0.6618'
$ws.Range("E45").Value = '  +2.22%  '

Set-TextValue $ws.Range("D46") '2.355'
$ws.Range("E46").Value = '  +4.22%  '

$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D47") '4.034'
$ws.Range("E47").Value = '  +0.61%  '

$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D48") '0.9962'
$ws.Range("E48").Value = '  -0.48%  '

$ws.Range("E49").Value = '  +1.45%  '

Set-TextValue $ws.Range("D50") '131.63'
$ws.Range("E50").Value = '  -0.12%  '

Set-TextValue $ws.Range("D51") '1.211'
$ws.Range("E51").Value = '  +0.61%  '

$excel.CutCopyMode = 0
